$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 currently holds the "culpeper / 52" record. Overwrite it in place with
# the data that used to live in row 4 ("Orange / 53"), then delete row 4 so
# the two rows collapse into one (row 3), matching the target layout.
$ws.Range("A3").Value = 53
$ws.Range("D3").Value = "127 Belleview Ave"
$ws.Range("E3").Value = "Orange"
$ws.Range("G3").Value = -78.113039999999998
$ws.Range("H3").Value = 38.245869999999996

# Remove the now-duplicated row 4 (its data has been folded into row 3 above).
$ws.Rows(4).Delete()

# Widen column D to fit the longer address text.
$ws.Columns("D").ColumnWidth = 19.43

# Move the active selection, matching the saved cursor position.
$ws.Range("C10").Select()
